$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E4").Value = "[5.00485605 5.98350954 6.31381953]"
$ws.Range("F4").Value = "[5.004856045049252, 5.983509544109174, 6.3138195277544105]"
$ws.Range("G4").Value = 380.107
$ws.Range("H4").Value = 380.107
$ws.Range("I4").Value = 19.4963
$ws.Range("J4").Value = 15.3183
$ws.Range("D7").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E7").Value = "[5.00099821 6.27957958 1.22416492]"
$ws.Range("F7").Value = "[5.000998205722637, 6.279579575394335, 1.2241649211971726]"
$ws.Range("G7").Value = 410.445
$ws.Range("H7").Value = 410.445
$ws.Range("I7").Value = 20.2594
$ws.Range("J7").Value = 15.831
$ws.Range("F8").Value = "[5.01816187654278, 6.283904691357364, 1.1355432098765736]"
$ws.Range("F9").Value = "[5.0166004938267905, 6.2858627160487455, 1.1116049382716284]"
$ws.Range("D10").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E10").Value = "[5.01610481 6.30204009 0.71287796]"
$ws.Range("F10").Value = "[5.016104807844141, 6.3020400947547195, 0.7128779631878824]"
$ws.Range("G10").Value = 15.8857
$ws.Range("H10").Value = 415.563
$ws.Range("I10").Value = 20.3854
$ws.Range("J10").Value = 15.8857
$ws.Range("K10").Value = 1881
$ws.Range("D13").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E13").Value = "[5.00629736 6.24306192 0.55978661]"
$ws.Range("F13").Value = "[5.0062973609303425, 6.243061921915732, 0.5597866051631579]"
$ws.Range("G13").Value = 15.9353
$ws.Range("H13").Value = 417.131
$ws.Range("I13").Value = 20.4238
$ws.Range("J13").Value = 15.9353
$ws.Range("K13").Value = 9690
$ws.Range("D16").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E16").Value = "[3.52956016 0.53776108 0.05254562]"
$ws.Range("F16").Value = "[3.529560159942636, 0.5377610824091528, 0.05254561750885019]"
$ws.Range("G16").Value = 88.9281
$ws.Range("H16").Value = 12280.6
$ws.Range("I16").Value = 110.818
$ws.Range("J16").Value = 88.9281
$ws.Range("F17").Value = "[5.001803905449378, 6.217746118086512, 2.2871260708284074]"
$ws.Range("F18").Value = "[5.001863826357645, 6.213147485618186, 2.3661798366712326]"
$ws.Range("D19").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E19").Value = "[5.00180884 6.21736774 2.29363062]"
$ws.Range("F19").Value = "[5.001808835744405, 6.217367742400251, 2.2936306182365147]"
$ws.Range("G19").Value = 20.0548
$ws.Range("H19").Value = 402.194
$ws.Range("I19").Value = 20.0548
$ws.Range("J19").Value = 15.6823
$ws.Range("F20").Value = "[5.000004796867698, 6.295228044651211, 0.5105980238884399]"
$ws.Range("F21").Value = "[5.000111977584725, 6.291770872525101, 0.6050154873708847]"
$ws.Range("D22").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E22").Value = "[5.00001128 6.2948113  0.51835745]"
$ws.Range("F22").Value = "[5.000011284357589, 6.29481130288869, 0.5183574481291469]"
$ws.Range("G22").Value = 20.4092
$ws.Range("H22").Value = 416.534
$ws.Range("I22").Value = 20.4092
$ws.Range("J22").Value = 15.9362
$ws.Range("F23").Value = "[4.451399478396651, 0.6352689964662742, 0.0367689962822413]"
$ws.Range("F24").Value = "[4.4955877486300215, 0.7459156620336217, 0.137504766871232]"
$ws.Range("D25").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E25").Value = "[4.45231125 0.63624176 0.04469071]"
$ws.Range("F25").Value = "[4.452311250944099, 0.636241762532293, 0.0446907062081126]"
$ws.Range("G25").Value = 80.43729999999999
$ws.Range("H25").Value = 6470.15
$ws.Range("I25").Value = 80.43729999999999
$ws.Range("J25").Value = 65.5408
$ws.Range("D28").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E28").Value = "[554.05666544  71.6327101   13.40511746]"
$ws.Range("F28").Value = "[4.983958202327891, 5.739800489022603, 636.1703007284555]"
$ws.Range("D31").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E31").Value = "[232.25702181  39.96976741  27.39333956]"
$ws.Range("F31").Value = "[2.089243503621952, 3.2027057223332913, 288.4206692870766]"
$ws.Range("G31").Value = 29037.6
$ws.Range("I31").Value = 319.485
$ws.Range("J31").Value = 276.28
$ws.Range("D34").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E34").Value = "[29.50347539  4.90865786  6.61596435]"
$ws.Range("F34").Value = "[0.2653953960281711, 0.3933219440985611, 39.774655972837294]"
$ws.Range("G34").Value = 77746.39999999999
$ws.Range("H34").Value = 76909.5
$ws.Range("I34").Value = 277.326
$ws.Range("D37").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E37").Value = "[3.77517193 0.0286407  3.63039648]"
$ws.Range("F37").Value = "[0.033959160272652826, 0.0022949280883051964, 7.875199718664754]"
$ws.Range("H37").Value = 85498.5
$ws.Range("I37").Value = 292.401
$ws.Range("D40").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E40").Value = "[ 0.38165351 -0.00599178  0.3906187 ]"
$ws.Range("F40").Value = "[0.003433123800215943, -0.0004801102636563552, 0.8197783833660495]"
$ws.Range("G40").Value = 216.416
$ws.Range("H40").Value = 86961.2
$ws.Range("I40").Value = 294.892
$ws.Range("F41").Value = "[0.0003412185028047634, -5.326661309567206e-05, 0.0809260484533906]"
$ws.Range("D43").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E43").Value = "[ 4.02602999e-02 -8.87010753e-06  4.61742595e-02]"
$ws.Range("F43").Value = "[0.00036215727497283286, -7.107457953594998e-07, 0.09144394728558917]"
$ws.Range("H43").Value = 87114.2
$ws.Range("I43").Value = 295.151
$ws.Range("J43").Value = 216.662
$ws.Range("D46").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E46").Value = "[5.14415988 0.84303779 1.2427115 ]"
$ws.Range("F46").Value = "[0.04627374672541609, 0.0675511048428729, 7.024227800986424]"
$ws.Range("G46").Value = 292.313
$ws.Range("H46").Value = 85056.89999999999
$ws.Range("I46").Value = 291.645
$ws.Range("D49").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E49").Value = "[0.51651862 0.08462336 0.13300942]"
$ws.Range("F49").Value = "[0.004646288700604582, 0.006780718036783926, 0.7135242767732808]"
$ws.Range("G49").Value = 294.893
$ws.Range("H49").Value = 86918.10000000001
$ws.Range("I49").Value = 294.819
$ws.Range("D52").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E52").Value = "[0.05374688 0.00904989 0.02042967]"
$ws.Range("F52").Value = "[0.00048347438765121447, 0.0007251516396662797, 0.08083496179066862]"
$ws.Range("H52").Value = 87109.89999999999
$ws.Range("I52").Value = 295.144
$ws.Range("D55").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E55").Value = "[288.09661897  48.82206751  36.98859438]"
$ws.Range("F55").Value = "[5.00783666633147, 5.733892675659728, 10.599070880235033]"
$ws.Range("D58").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E58").Value = "[249.1284227   42.21838806  31.98651382]"
$ws.Range("F58").Value = "[4.3304723751343115, 4.958325577723893, 9.16641991352505]"
$ws.Range("G58").Value = 1960.5
$ws.Range("H58").Value = 1960.5
$ws.Range("I58").Value = 44.2775
$ws.Range("J58").Value = 33.7467
$ws.Range("D61").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E61").Value = "[52.22822875  8.85102985  6.711784  ]"
$ws.Range("F61").Value = "[0.9078566762123806, 1.0395064733523458, 1.9274589091423904]"
$ws.Range("G61").Value = 58528
$ws.Range("H61").Value = 58528
$ws.Range("I61").Value = 241.926
$ws.Range("D64").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E64").Value = "[6.78537887 1.26521961 3.63550759]"
$ws.Range("F64").Value = "[0.11794678196789525, 0.14859332689740506, 2.883384959917618]"
$ws.Range("H64").Value = 82891.2
$ws.Range("I64").Value = 287.908
$ws.Range("D67").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E67").Value = "[0.68168179 0.12036249 0.3906187 ]"
$ws.Range("F67").Value = "[0.011849327113609684, 0.014135935183670543, 0.3226952075492269]"
$ws.Range("H67").Value = 86700
$ws.Range("I67").Value = 294.449
$ws.Range("D70").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E70").Value = "[0.07026313 0.01262656 0.04617426]"
$ws.Range("F70").Value = "[0.0012213481467994993, 0.0014829220009943606, 0.03892362750530713]"
$ws.Range("H70").Value = 87087.7
$ws.Range("I70").Value = 295.106
$ws.Range("D73").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E73").Value = "[9.76209664 1.65457862 1.26069737]"
$ws.Range("F73").Value = "[0.16968954960475457, 0.1943214767194512, 0.3662076573839981]"
$ws.Range("G73").Value = 285.223
$ws.Range("H73").Value = 81352.39999999999
$ws.Range("I73").Value = 285.223
$ws.Range("D76").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E76").Value = "[0.97837059 0.16605935 0.13319059]"
$ws.Range("F76").Value = "[0.017006517202371917, 0.01950278922254529, 0.043277298027593336]"
$ws.Range("G76").Value = 294.184
$ws.Range("H76").Value = 86544
$ws.Range("I76").Value = 294.184
$ws.Range("D79").Value = "[0.002327721424381168, 0.0006558972239079049, 0.007902654545887318]"
$ws.Range("E79").Value = "[0.09993265 0.01719635 0.02043153]"
$ws.Range("F79").Value = "[0.0017370782412306723, 0.002019620097596851, 0.0109818599995041]"
$ws.Range("G79").Value = 295.08
$ws.Range("H79").Value = 87072
$ws.Range("I79").Value = 295.08
